$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the last existing row (308): high/close were revised ---
$ws.Range("D308").Value = 7.808
$ws.Range("F308").Value = 7.808

# --- Give the three new rows (309:311) the same look as the data above them
#     (date column A uses style index 2: centered, bordered, datetime format) ---
$ws.Range("A308").Copy()
$ws.Range("A309:A311").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 309 ---
$ws.Range("A309").Value = 45047.33333333334
$ws.Range("B309").Value = "FX_IDC:USDGTQ"
$ws.Range("C309").Value = 7.809
$ws.Range("D309").Value = 7.83
$ws.Range("E309").Value = 7.8015
$ws.Range("F309").Value = 7.83
$ws.Range("G309").Value = 0

# --- Row 310 ---
$ws.Range("A310").Value = 45078.33333333334
$ws.Range("B310").Value = "FX_IDC:USDGTQ"
$ws.Range("C310").Value = 7.83
$ws.Range("D310").Value = 7.8455
$ws.Range("E310").Value = 7.829
$ws.Range("F310").Value = 7.8455
$ws.Range("G310").Value = 0

# --- Row 311 ---
$ws.Range("A311").Value = 45110.33333333334
$ws.Range("B311").Value = "FX_IDC:USDGTQ"
$ws.Range("C311").Value = 7.8455
$ws.Range("D311").Value = 7.847
$ws.Range("E311").Value = 7.845
$ws.Range("F311").Value = 7.847
$ws.Range("G311").Value = 0
